$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# Position / size the subtitle placeholder explicitly (was inherited, now fixed)
$shp.Left = 120
$shp.Top = 425.2705
$shp.Width = 720
$shp.Height = 29.625

# Turn off autofit (normAutofit -> noAutofit)
$shp.TextFrame.AutoSize = 0

# Add the "Sparks Lu" author line above the existing "Last updated" line,
# then resize both lines down to 14pt (from 17pt).
$tr = $shp.TextFrame.TextRange
$tr.InsertBefore("Sparks Lu`r")
$tr.Font.Size = 14
